$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 - shape "직사각형 8" (id=9)
#   "...장르 선호의 검정 결과의 " -> "...장르 선호의 " + "카이제곱검정" + " "
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(1)
$tr6 = $shp6.TextFrame.TextRange

# split "검정 결과의" (6 chars starting at position 23) into its own run with new text
$tr6.Characters(23, 6).Text = "카이제곱검정"
# split the following space (position 29) into its own run (keeps same character)
$tr6.Characters(29, 1).Text = " "

# restore the autosized shape back to its original geometry (unaffected by this
# edit in the target deck - only slide 12's box changes size)
$shp6.Left = 63.68818897637795
$shp6.Top = 187.72283464566928
$shp6.Width = 802.9785039370079
$shp6.Height = 159.9468503937008

# ---------------------------------------------------------------------------
# Slide 12 - shape "직사각형 8" (id=9)
#   "...선호(트렌드)의 검정 결과의 " -> "...선호(트렌드)의 " + "카이제곱검정" + " "
#   Also drop the trailing empty paragraph and the
#   "지역별로도 연도별 트렌드에 차이가 있는지 분석함" paragraph, and shrink the
#   shape to its new autofit height.
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(1)
$tr12 = $shp12.TextFrame.TextRange

$tr12.Characters(28, 6).Text = "카이제곱검정"
$tr12.Characters(34, 1).Text = " "

# remove the last two paragraphs: the blank paragraph and the
# "지역별로도 연도별 트렌드에 차이가 있는지 분석함" paragraph (28 chars incl. the
# two leading paragraph marks)
$tr12.Characters(173, 28).Text = ""

# resize the shape to match the now-shorter text block
$shp12.Left = 61.996929133858266
$shp12.Top = 166.38953455905514
$shp12.Width = 802.9785039370079
$shp12.Height = 138.13590551181102
